# Auto-generated Excel COM-interop script applying the scheduled-runner price update.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the refreshed market data.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 967.4375
$ws.Range("J80").Value = 999.5
$ws.Range("L80").Value = 2998.5
$ws.Range("N80").Value = -4994.5
$ws.Range("H83").Value = 967.4375
$ws.Range("J83").Value = 999.5
$ws.Range("L83").Value = 8995.5
$ws.Range("N83").Value = -18979.5
$ws.Range("H94").Value = 5554.375
$ws.Range("I94").Value = 2062.1428
$ws.Range("K94").Value = 2062.1428
$ws.Range("M94").Value = -1611.1428
$ws.Range("H99").Value = 534.4
$ws.Range("J99").Value = 214.5
$ws.Range("L99").Value = 643.5
$ws.Range("N99").Value = -3639.5
$ws.Range("H107").Value = 1210.871
$ws.Range("I107").Value = 944.1539
$ws.Range("J107").Value = 2597.8
$ws.Range("K107").Value = 944.1539
$ws.Range("L107").Value = 2597.8
$ws.Range("M107").Value = 975.8461
$ws.Range("N107").Value = -6437.8
$ws.Range("H135").Value = 1919.6774
$ws.Range("I135").Value = 1622.1786
$ws.Range("K135").Value = 14599.6074
$ws.Range("M135").Value = -12064.6074

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2232.5557
$ws.Range("I45").Value = 1096.2941
$ws.Range("K45").Value = 1096.2941
$ws.Range("M45").Value = -719.2941000000001
$ws.Range("H61").Value = 7420.722
$ws.Range("I61").Value = 1161.1818
$ws.Range("K61").Value = 1161.1818
$ws.Range("M61").Value = -949.1818000000001
$ws.Range("H102").Value = 2002.1364
$ws.Range("I102").Value = 1887.4286
$ws.Range("K102").Value = 1887.4286
$ws.Range("M102").Value = -265.4286
$ws.Range("H110").Value = 3507.6511
$ws.Range("I110").Value = 3421.7693
$ws.Range("K110").Value = 3421.7693
$ws.Range("M110").Value = -1376.7693
$ws.Range("H111").Value = 58214.332
$ws.Range("J111").Value = 58214.332
$ws.Range("L111").Value = 58214.332
$ws.Range("N111").Value = -66394.33199999999
$ws.Range("H122").Value = 1456.4286
$ws.Range("I122").Value = 1240.4584
$ws.Range("K122").Value = 3721.3752
$ws.Range("M122").Value = -1271.3752
$ws.Range("H132").Value = 1260.0513
$ws.Range("I132").Value = 677.0606
$ws.Range("K132").Value = 2031.1818
$ws.Range("M132").Value = 498.8181999999999
$ws.Range("H136").Value = 7420.722
$ws.Range("I136").Value = 1161.1818
$ws.Range("K136").Value = 3483.5454
$ws.Range("M136").Value = -933.5454

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1969.262
$ws.Range("I107").Value = 1910.2703
$ws.Range("K107").Value = 1910.2703
$ws.Range("M107").Value = 9.729700000000093

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 5051
$ws.Range("I39").Value = 5051
$ws.Range("K39").Value = 5051
$ws.Range("M39").Value = -4660
$ws.Range("H49").Value = 5051
$ws.Range("I49").Value = 5051
$ws.Range("K49").Value = 5051
$ws.Range("M49").Value = -4869
$ws.Range("H86").Value = 84542.78
$ws.Range("I86").Value = 130577.4
$ws.Range("K86").Value = 130577.4
$ws.Range("M86").Value = -129454.4
$ws.Range("H89").Value = 84542.78
$ws.Range("I89").Value = 130577.4
$ws.Range("K89").Value = 652887
$ws.Range("M89").Value = -647271
$ws.Range("H94").Value = 2098.5557
$ws.Range("J94").Value = 2788.2222
$ws.Range("L94").Value = 2788.2222
$ws.Range("N94").Value = -3690.2222
$ws.Range("H105").Value = 1475.6875
$ws.Range("I105").Value = 972.2857
$ws.Range("K105").Value = 972.2857
$ws.Range("M105").Value = 774.7143
$ws.Range("H132").Value = 51944.45
$ws.Range("I132").Value = 72277.07000000001
$ws.Range("J132").Value = 4501.6665
$ws.Range("K132").Value = 216831.21
$ws.Range("L132").Value = 13504.9995
$ws.Range("M132").Value = -214301.21
$ws.Range("N132").Value = -18564.9995
$ws.Range("H134").Value = 3143.1365
$ws.Range("I134").Value = 2807.842
$ws.Range("J134").Value = 5266.6665
$ws.Range("K134").Value = 8423.526
$ws.Range("L134").Value = 15799.9995
$ws.Range("M134").Value = -5888.526
$ws.Range("N134").Value = -20869.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 12197.255
$ws.Range("I99").Value = 3425
$ws.Range("K99").Value = 10275
$ws.Range("M99").Value = -8029
$ws.Range("H113").Value = 1186.875
$ws.Range("J113").Value = 1213.5714
$ws.Range("L113").Value = 3640.7142
$ws.Range("N113").Value = -7980.7142
$ws.Range("H129").Value = 2438.4
$ws.Range("J129").Value = 2243.2307
$ws.Range("L129").Value = 6729.6921
$ws.Range("N129").Value = -16729.6921
$ws.Range("H132").Value = 1916.8182
$ws.Range("J132").Value = 2273.25
$ws.Range("L132").Value = 20459.25
$ws.Range("N132").Value = -25519.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5144.909
$ws.Range("I70").Value = 5957.4
$ws.Range("J70").Value = 4467.8335
$ws.Range("K70").Value = 5957.4
$ws.Range("L70").Value = 4467.8335
$ws.Range("M70").Value = -5687.4
$ws.Range("N70").Value = -5007.8335
$ws.Range("H73").Value = 5144.909
$ws.Range("I73").Value = 5957.4
$ws.Range("J73").Value = 4467.8335
$ws.Range("K73").Value = 5957.4
$ws.Range("L73").Value = 4467.8335
$ws.Range("M73").Value = -5021.4
$ws.Range("N73").Value = -6339.8335
$ws.Range("H113").Value = 2998.6428
$ws.Range("I113").Value = 3248.375
$ws.Range("J113").Value = 2665.6667
$ws.Range("K113").Value = 3248.375
$ws.Range("L113").Value = 2665.6667
$ws.Range("M113").Value = -1078.375
$ws.Range("N113").Value = -7005.6667
$ws.Range("H122").Value = 3410.1667
$ws.Range("I122").Value = 3098.9565
$ws.Range("J122").Value = 4432.7144
$ws.Range("K122").Value = 9296.869499999999
$ws.Range("L122").Value = 13298.1432
$ws.Range("M122").Value = -6846.869499999999
$ws.Range("N122").Value = -18198.1432
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()  # was 49069
$ws.Range("N124").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()  # was 100000
$ws.Range("N127").Value = 0
$ws.Range("H132").Value = 2401.7693
$ws.Range("I132").Value = 1722.4445
$ws.Range("K132").Value = 5167.333500000001
$ws.Range("M132").Value = -2637.333500000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2499.6667
$ws.Range("I40").Value = 2499.6667
$ws.Range("K40").Value = 2499.6667
$ws.Range("M40").Value = -2363.6667
$ws.Range("H46").Value = 2561.2083
$ws.Range("I46").Value = 945.1111
$ws.Range("J46").Value = 3530.8667
$ws.Range("K46").Value = 945.1111
$ws.Range("L46").Value = 3530.8667
$ws.Range("M46").Value = -757.1111
$ws.Range("N46").Value = -3906.8667
$ws.Range("H63").Value = 99999
$ws.Range("J63").Value = 99999
$ws.Range("L63").Value = 99999
$ws.Range("N63").Value = -101497
$ws.Range("H66").Value = 99999
$ws.Range("J66").Value = 99999
$ws.Range("L66").Value = 299997
$ws.Range("N66").Value = -307485
$ws.Range("H100").Value = 2508.739
$ws.Range("J100").Value = 2978.6
$ws.Range("L100").Value = 2978.6
$ws.Range("N100").Value = -4060.6
$ws.Range("H122").Value = 3051
$ws.Range("I122").Value = 3017.8333
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 9053.499899999999
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -6603.499899999999
$ws.Range("N122").Value = -14650

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 953.875
$ws.Range("I113").Value = 807.0714
$ws.Range("K113").Value = 2421.2142
$ws.Range("M113").Value = -251.2142000000003
$ws.Range("H122").Value = 10028563
$ws.Range("I122").Value = 11142582
$ws.Range("J122").Value = 2396
$ws.Range("K122").Value = 33427746
$ws.Range("L122").Value = 7188
$ws.Range("M122").Value = -33425296
$ws.Range("N122").Value = -12088
$ws.Range("H125").Value = 53742.5
$ws.Range("J125").Value = 53742.5
$ws.Range("L125").Value = 53742.5
$ws.Range("N125").Value = -63582.5
$ws.Range("H132").Value = 20892866
$ws.Range("I132").Value = 22792152
$ws.Range("K132").Value = 68376456
$ws.Range("M132").Value = -68373926
